$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header for column H, matching the style of G1 (antecedents_length header)
$ws.Range("H1").Value = "consequents_length"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill H2:H34 with the value 1 (consequents_length for every rule)
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
